$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate so no formatting bleeds over from the old layout
$ws.Cells.Clear()

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "query"
$ws.Range("B1").Value = "TabName"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# --- Row 2 (CasesTab) ---
$ws.Range("B2").Value = "CasesTab"
$ws.Range("E2").Value = "TC03_Trials_Filter_AssocFileFormat-Vcf_WebData.xlsx"
$caseQuery = @"
MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s WHERE f.file_format IN ['vcf']  RETURN DISTINCT coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity
"@
$ws.Range("A2").Value = $caseQuery
$casesTabQuery = @"
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
WHERE f.file_format IN ['vcf'] 
RETURN DISTINCT
    c.case_id AS ``Case ID``,
     ct.clinical_trial_designation AS ``Trial Code``,
     a.arm_id AS Arm,
      a.arm_drug AS ``Arm Treatment``,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
"@
$ws.Range("C2").Value = $casesTabQuery
$statQuery = @"
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE f.file_format IN ['vcf'] 
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
"@
$ws.Range("D2").Value = $statQuery

# --- Row 3 (FilesTab) ---
$ws.Range("B3").Value = "FilesTab"
$ws.Range("E3").Value = "TC03_Trials_Filter_AssocFileFormat-Vcf_WebData.xlsx"
$filesTabQuery = @"
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
 WHERE f.file_format IN ['vcf'] 
 WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS ``File Name``,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS ``File Format``,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS ``Trial Code``,
    a.arm_id AS Arm,
    c.case_id AS ``Case ID``
"@
$ws.Range("C3").Value = $filesTabQuery
$ws.Range("D3").Value = $statQuery

# --- Formatting: wrap text on the long-text cells ---
$ws.Range("A2").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Range("D2").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("D3").WrapText = $true

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 75.81640625
$ws.Columns.Item(2).ColumnWidth = 8.81640625
$ws.Columns.Item(3).ColumnWidth = 75.81640625
$ws.Columns.Item(4).ColumnWidth = 70.26953125
$ws.Columns.Item(5).ColumnWidth = 28.54296875

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 188.5
$ws.Rows.Item(3).RowHeight = 409.5

# --- View state ---
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("D2").Select()

Write-Host "edit applied"
